{"js": "// Remove the trailing \"Ver no Jupiter...\" paragraph, the \"\u00a9 2020...\"\n// paragraph, and the blank paragraph that immediately precedes them\n// (right after the \"LOQ4205: ... (Requisito fraco)\" paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4205\" requirement paragraph and the two text\n// paragraphs that must be removed, by their exact text content.\nlet loqIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"LOQ4205\") !== -1) {\n    loqIndex = i;\n  } else if (t === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiterIndex = i;\n  } else if (t.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nconst toDelete = [];\nif (jupiterIndex !== -1) toDelete.push(jupiterIndex);\nif (copyrightIndex !== -1) toDelete.push(copyrightIndex);\n// The blank paragraph directly after the LOQ4205 paragraph (and right\n// before the \"Ver no Jupiter...\" paragraph) is removed too.\nif (loqIndex !== -1 && loqIndex + 1 < items.length && items[loqIndex + 1].text === \"\") {\n  toDelete.push(loqIndex + 1);\n}\n\n// Delete from the highest index down so earlier indices stay valid.\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" paragraph, the \"(c) 2020...\"\n# paragraph, and the blank paragraph that immediately precedes them\n# (right after the \"LOQ4205: ... (Requisito fraco)\" paragraph).\n$d = $word.ActiveDocument\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n$loqIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    $trimmed = $t.Trim()\n    if ($trimmed -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $jupiterIndex = $i\n    } elseif ($t -like \"*Powered by Jekyll and Github pages*\") {\n        $copyrightIndex = $i\n    } elseif ($t -like \"*LOQ4205*\") {\n        $loqIndex = $i\n    }\n}\n\n$toDelete = @()\nif ($jupiterIndex -ne -1) { $toDelete += $jupiterIndex }\nif ($copyrightIndex -ne -1) { $toDelete += $copyrightIndex }\nif ($loqIndex -ne -1) {\n    $nextIndex = $loqIndex + 1\n    if ($nextIndex -le $count -and $d.Paragraphs.Item($nextIndex).Range.Text.Trim() -eq \"\") {\n        $toDelete += $nextIndex\n    }\n}\n\n# Delete from the highest paragraph index down so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete() | Out-Null\n}\n"}
